$wb = $excel.ActiveWorkbook

# --- Sheet: Means ---
$wsMeans = $wb.Worksheets.Item("Means")

# Row 9: Total Cancer Risk (per million)
$wsMeans.Range("B9").Value = 23
$wsMeans.Range("C9").Value = 26
$wsMeans.Range("G9").Value = 30

# Row 10: Total Respiratory (hazard quotient)
$wsMeans.Range("B10").Value = 0.27
$wsMeans.Range("C10").Value = 0.28
$wsMeans.Range("D10").Value = 0.3
$wsMeans.Range("E10").Value = 0.3
$wsMeans.Range("F10").Value = 0.31
$wsMeans.Range("G10").Value = 0.34

# --- Sheet: Standard Deviations ---
$wsSD = $wb.Worksheets.Item("Standard Deviations")

# Row 9: Total Cancer Risk (per million) SD
$wsSD.Range("B9").Value = 7.2
$wsSD.Range("C9").Value = 6.3
$wsSD.Range("G9").Value = 3.6

# Row 10: Total Respiratory (hazard quotient) SD
$wsSD.Range("B10").Value = 0.094
$wsSD.Range("C10").Value = 0.058
$wsSD.Range("E10").Value = 0
$wsSD.Range("F10").Value = 0.032
$wsSD.Range("G10").Value = 0.05
